$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 154; this shifts the existing rows 154:248
# down to 155:249 and copies formatting from the row above (row 153),
# which gives the date column (D) the correct date number format.
$ws.Rows("154:154").Insert()

# Populate the newly inserted row 154 with the new weekly record.
$ws.Range("A154").Value2 = 3
$ws.Range("B154").Value2 = "Femacal de La Calera"
$ws.Range("C154").Value2 = "Coquimbo"
$ws.Range("D154").Value2 = 44875
$ws.Range("E154").Value2 = 5
$ws.Range("F154").Value2 = "Fruta"
$ws.Range("G154").Value2 = 100101
$ws.Range("H154").Value2 = "Berries"
$ws.Range("I154").Value2 = 100101001
$ws.Range("J154").Value2 = "Arándano (blue)"
$ws.Range("K154").Value2 = "Sin especificar"
$ws.Range("L154").Value2 = "Primera"
$ws.Range("M154").Value2 = 57
$ws.Range("N154").Value2 = 8000
$ws.Range("O154").Value2 = 8000
$ws.Range("P154").Value2 = 8000
$ws.Range("Q154").Value2 = "$/bandeja 2 kilos"
$ws.Range("R154").Value2 = "Provincia de Quillota"
$ws.Range("S154").Value2 = 4000
$ws.Range("T154").Value2 = 2
